$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the original "_GoBack" bookmark near "Dec 2, 2019".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Merge the 3 runs "Appen"/"d"/"ix 3" (inside the hyperlink pointing to
#    Appendix 3) into a single run, preserving the Hyperlink character style.
#    (A two-step replace forces the run-merge even though the final text is
#    identical to the concatenation of the original runs.)
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Appendix 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Text = "ZZZ_PLACEHOLDER_ZZZ"
$rng2b = $d.Content
$rng2b.Find.Execute("ZZZ_PLACEHOLDER_ZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2b.Text = "Appendix 3"
$rng2b.Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 3) Move the "UploadedFiles: a folder containing:" / "RSIND: a folder
#    holding..." bullets from the numId=4 (outside-web-tree) list up into the
#    numId=5 (web-tree) list, right after "RSINDUpload.php" and before
#    "jqUpload". First delete the two paragraphs from their old spot.
# ---------------------------------------------------------------------------
$rngDel1 = $d.Content
$rngDel1.Find.Execute("UploadedFiles: a folder containing:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngDel1.Paragraphs(1).Range.Delete()

$rngDel2 = $d.Content
$rngDel2.Find.Execute("RSIND: a folder holding the RSIND file uploaded by", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngDel2.Paragraphs(1).Range.Delete()

# Now insert the two replacement paragraphs right before the "jqUpload" bullet
# (i.e. at the end of the preceding "RSINDUpload.php..." paragraph) using the
# same numId=5 list, one level 0 and one level 1, each carrying the same
# spell-check proofErr markers Word itself would emit.
$rngIns = $d.Content
$rngIns.Find.Execute("drop area are shown above.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($rngIns.End, $rngIns.End)

$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:pStyle w:val="Normal1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>UploadedFiles</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: a folder containing:</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Normal1"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">RSIND: a folder holding the RSIND files uploaded by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RSINDUpload.php</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">.  It’s held there until validated and (possibly) copied by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TestInstallRSIND.bash</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.  Files stored in this directory remain until they are “old”, thus preventing the accidental upload of the same file twice.</w:t></w:r></w:p>
</w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insertPoint.InsertXML($xmlSnippet)
